$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 5057.75
$ws.Range("J7").Value = 5743.6665
$ws.Range("L7").Value = 5743.6665
$ws.Range("N7").Value = -5967.6665
$ws.Range("H14").Value = 5057.75
$ws.Range("J14").Value = 5743.6665
$ws.Range("L14").Value = 5743.6665
$ws.Range("N14").Value = -6125.6665
$ws.Range("H16").Value = 1052.25
$ws.Range("J16").Value = 1500
$ws.Range("L16").Value = 1500
$ws.Range("N16").Value = -1960
$ws.Range("H58").Value = 1041
$ws.Range("I58").Value = 109
$ws.Range("J58").Value = 1600.2
$ws.Range("K58").Value = 327
$ws.Range("L58").Value = 4800.6
$ws.Range("M58").Value = -177
$ws.Range("N58").Value = -5100.6
$ws.Range("H62").Value = 10908.363
$ws.Range("I62").Value = 8498
$ws.Range("K62").Value = 8498
$ws.Range("M62").Value = -7874
$ws.Range("H65").Value = 10908.363
$ws.Range("I65").Value = 8498
$ws.Range("K65").Value = 42490
$ws.Range("M65").Value = -39370
$ws.Range("H70").Value = 5434.25
$ws.Range("I70").Value = 1737.5
$ws.Range("J70").Value = 6666.5
$ws.Range("K70").Value = 5212.5
$ws.Range("L70").Value = 19999.5
$ws.Range("M70").Value = -4942.5
$ws.Range("N70").Value = -20539.5
$ws.Range("H73").Value = 5434.25
$ws.Range("I73").Value = 1737.5
$ws.Range("J73").Value = 6666.5
$ws.Range("K73").Value = 5212.5
$ws.Range("L73").Value = 19999.5
$ws.Range("M73").Value = -4276.5
$ws.Range("N73").Value = -21871.5
$ws.Range("H96").Value = 307.66666
$ws.Range("I96").Value = 204.2
$ws.Range("J96").Value = 825
$ws.Range("K96").Value = 612.5999999999999
$ws.Range("L96").Value = 2475
$ws.Range("M96").Value = 760.4000000000001
$ws.Range("N96").Value = -5221
$ws.Range("H98").Value = 1592.8823
$ws.Range("I98").Value = 1592.8823
$ws.Range("K98").Value = 1592.8823
$ws.Range("M98").Value = -94.88229999999999
$ws.Range("H112").Value = 2976.8333
$ws.Range("I112").Value = 1997.5
$ws.Range("J112").Value = 3172.7
$ws.Range("K112").Value = 5992.5
$ws.Range("L112").Value = 9518.099999999999
$ws.Range("M112").Value = -4884.5
$ws.Range("N112").Value = -11734.1
$ws.Range("H122").Value = 1592.8823
$ws.Range("I122").Value = 1592.8823
$ws.Range("K122").Value = 4778.6469
$ws.Range("M122").Value = -2328.6469
$ws.Range("H131").Value = 1379.875
$ws.Range("I131").Value = 1173.3334
$ws.Range("J131").Value = 1999.5
$ws.Range("K131").Value = 3520.0002
$ws.Range("L131").Value = 5998.5
$ws.Range("M131").Value = 1519.9998
$ws.Range("N131").Value = -16078.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 10046
$ws.Range("I2").Value = 5670.8
$ws.Range("J2").Value = 17338
$ws.Range("K2").Value = 5670.8
$ws.Range("L2").Value = 17338
$ws.Range("M2").Value = -5557.8
$ws.Range("N2").Value = -17564
$ws.Range("H9").Value = 10500
$ws.Range("I9").Value = 7500
$ws.Range("J9").Value = 12000
$ws.Range("K9").Value = 7500
$ws.Range("L9").Value = 12000
$ws.Range("M9").Value = -7330
$ws.Range("N9").Value = -12340
$ws.Range("H20").Value = 10500
$ws.Range("I20").Value = 7500
$ws.Range("J20").Value = 12000
$ws.Range("K20").Value = 7500
$ws.Range("L20").Value = 12000
$ws.Range("M20").Value = -7230
$ws.Range("N20").Value = -12540
$ws.Range("H23").Value = 55006
$ws.Range("J23").Value = 55006
$ws.Range("L23").Value = 55006
$ws.Range("N23").Value = -55524
$ws.Range("H116").Value = 10046
$ws.Range("I116").Value = 5670.8
$ws.Range("J116").Value = 17338
$ws.Range("K116").Value = 5670.8
$ws.Range("L116").Value = 17338
$ws.Range("M116").Value = -3376.8
$ws.Range("N116").Value = -21926

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 10046
$ws.Range("I3").Value = 5670.8
$ws.Range("J3").Value = 17338
$ws.Range("K3").Value = 5670.8
$ws.Range("L3").Value = 17338
$ws.Range("M3").Value = -5556.8
$ws.Range("N3").Value = -17566
$ws.Range("H5").Value = 9673.143
$ws.Range("I5").Value = 106.4
$ws.Range("J5").Value = 33590
$ws.Range("K5").Value = 106.4
$ws.Range("L5").Value = 33590
$ws.Range("M5").Value = 6.599999999999994
$ws.Range("N5").Value = -33816
$ws.Range("H134").Value = 2855.8096
$ws.Range("I134").Value = 2442.9443
$ws.Range("K134").Value = 7328.8329
$ws.Range("M134").Value = -4793.8329

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6286.517
$ws.Range("I31").Value = 2958.3
$ws.Range("J31").Value = 8038.2104
$ws.Range("K31").Value = 2958.3
$ws.Range("L31").Value = 8038.2104
$ws.Range("M31").Value = -2663.3
$ws.Range("N31").Value = -8628.2104
$ws.Range("H34").Value = 6286.517
$ws.Range("I34").Value = 2958.3
$ws.Range("J34").Value = 8038.2104
$ws.Range("K34").Value = 2958.3
$ws.Range("L34").Value = 8038.2104
$ws.Range("M34").Value = -2756.3
$ws.Range("N34").Value = -8442.2104
$ws.Range("H99").Value = 3907.7368
$ws.Range("I99").Value = 3759.2
$ws.Range("J99").Value = 4464.75
$ws.Range("K99").Value = 3759.2
$ws.Range("L99").Value = 4464.75
$ws.Range("M99").Value = -2261.2
$ws.Range("N99").Value = -7460.75
$ws.Range("H105").Value = 957.3
$ws.Range("I105").Value = 954.44446
$ws.Range("K105").Value = 954.44446
$ws.Range("M105").Value = 792.55554
$ws.Range("H126").Value = 3907.7368
$ws.Range("I126").Value = 3759.2
$ws.Range("J126").Value = 4464.75
$ws.Range("K126").Value = 11277.6
$ws.Range("L126").Value = 13394.25
$ws.Range("M126").Value = -8807.599999999999
$ws.Range("N126").Value = -18334.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5002.3125
$ws.Range("I70").Value = 4503.3335
$ws.Range("K70").Value = 4503.3335
$ws.Range("M70").Value = -4233.3335
$ws.Range("H73").Value = 5002.3125
$ws.Range("I73").Value = 4503.3335
$ws.Range("K73").Value = 4503.3335
$ws.Range("M73").Value = -3567.3335
$ws.Range("H86").Value = 45000
$ws.Range("J86").Value = 45000
$ws.Range("L86").Value = 45000
$ws.Range("N86").Value = -47372
$ws.Range("H89").Value = 45000
$ws.Range("J89").Value = 45000
$ws.Range("L89").Value = 135000
$ws.Range("N89").Value = -146856
$ws.Range("H122").Value = 360540
$ws.Range("I122").Value = 457299.38
$ws.Range("K122").Value = 1371898.14
$ws.Range("M122").Value = -1369448.14

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 3332.6667
$ws.Range("I10").Value = 3499
$ws.Range("J10").Value = 3000
$ws.Range("K10").Value = 3499
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = -3359
$ws.Range("N10").Value = -3280
$ws.Range("H40").Value = 7123.875
$ws.Range("I40").Value = 5399.8
$ws.Range("K40").Value = 5399.8
$ws.Range("M40").Value = -5263.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 420.75
$ws.Range("J6").Value = 675
$ws.Range("L6").Value = 675
$ws.Range("N6").Value = -905
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
